# Update the maturity assessment values and add explanatory notes.
$wb = $excel.ActiveWorkbook
$wsAval  = $wb.Worksheets.Item("Avaliação")
$wsResumo = $wb.Worksheets.Item("Resumo")

# Score updates on the "Avaliação" sheet
$wsAval.Range("C3").Value = 2
$wsAval.Range("C4").Value = 1

# Observações (column D) notes for each capability row
$wsAval.Range("D2").Value = "✅ Processo consolidado e amplamente adotado pelas equipes. Revisão periódica garante conformidade."
$wsAval.Range("D3").Value = "🔄 Ferramentas ativas, mas ainda com uso limitado por áreas não técnicas. Treinamentos adicionais recomendados."
$wsAval.Range("D4").Value = "🚧 Processo ainda não implementado. Requer definição clara de centros de custo e critérios de rateio."
$wsAval.Range("D5").Value = "✅ Orçamentos previstos regularmente com base em tendências históricas. Forecast validado com stakeholders."
$wsAval.Range("D6").Value = "🔄 Estratégias aplicadas em parte da infraestrutura. Potencial de expansão com análise mais granular de workloads."
$wsAval.Range("D7").Value = "🚧 Não há rotina estabelecida para reavaliação de workloads. Sugere-se calendário de revisões trimestrais."
$wsAval.Range("D8").Value = "🔄 Reuniões realizadas com certa frequência, mas ainda com baixa participação intersetorial."
$wsAval.Range("D9").Value = "🚧 Ausência de visualizações centralizadas e acessíveis. Recomendado criar painel com métricas críticas (tags, budget, anomalias)."

# Widen column D so the new notes are readable
$wsAval.Columns.Item(4).ColumnWidth = 98.88671875

# Move the active selection to D16 on the Avaliação sheet
$wsAval.Activate()
$wsAval.Range("D16").Select()
